# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-column suffixes to "_FV2404" / "_FV2410"
# - Freeze the header row
# - Wrap the data range in an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row strings -----------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2404"
    "B1" = "Segmentgruppe_FV2404"
    "C1" = "Segment_FV2404"
    "D1" = "Datenelement_FV2404"
    "E1" = "Segment ID_FV2404"
    "F1" = "Code_FV2404"
    "G1" = "Qualifier_FV2404"
    "H1" = "Beschreibung_FV2404"
    "I1" = "Bedingungsausdruck_FV2404"
    "J1" = "Bedingung_FV2404"
    "L1" = "Segmentname_FV2410"
    "M1" = "Segmentgruppe_FV2410"
    "N1" = "Segment_FV2410"
    "O1" = "Datenelement_FV2410"
    "P1" = "Segment ID_FV2410"
    "Q1" = "Code_FV2410"
    "R1" = "Qualifier_FV2410"
    "S1" = "Beschreibung_FV2410"
    "T1" = "Bedingungsausdruck_FV2410"
    "U1" = "Bedingung_FV2410"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into an Excel Table -----------------------------
$dataRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

Write-Host "Renamed headers, froze header row, and created Table1 over $($tbl.Range.Address())"
